$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.276005029678345
$ws.Range("B1").Value = 2.687077760696411
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.674093961715698
$ws.Range("E1").Value = 1.119034171104431
